$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "F3" = 2.730995758697492
    "G3" = 1.722369395643906
    "H3" = 2.878921494922184
    "I3" = 1.859543089110336
    "B4" = 65.292
    "C4" = 60.727
    "D4" = 56.775
    "E4" = 55
    "F4" = 114.321
    "G4" = 82.09
    "H4" = 118.689
    "I4" = 86.77800000000001
    "J4" = 53.345
    "K4" = 50.355
    "L4" = 47.736
    "M4" = 43.382
    "N4" = 39.921
    "O4" = 37.108
    "B5" = 32.032
    "C5" = 29.793
    "D5" = 27.854
    "E5" = 26.983
    "F5" = 56.086
    "G5" = 40.273
    "H5" = 58.229
    "I5" = 42.573
    "J5" = 26.171
    "K5" = 24.704
    "L5" = 23.419
    "M5" = 21.283
    "N5" = 19.585
    "O5" = 18.205
    "B6" = 1990.118
    "C6" = 1850.976
    "D6" = 1730.518
    "E6" = 1676.415
    "F6" = 1275.921
    "G6" = 1452.723
    "H6" = 1256.607
    "I6" = 1422.402
    "J6" = 1625.97
    "K6" = 1534.834
    "L6" = 1455.006
    "M6" = 1322.295
    "N6" = 1216.803
    "O6" = 1131.062
    "B7" = 26.642
    "C7" = 24.78
    "D7" = 23.167
    "E7" = 22.443
    "F7" = 46.649
    "G7" = 33.497
    "H7" = 48.431
    "I7" = 35.41
    "J7" = 21.767
    "K7" = 20.547
    "L7" = 19.478
    "M7" = 17.702
    "N7" = 16.29
    "O7" = 15.142
    "B8" = 4.816
    "C8" = 5.178
    "D8" = 5.539
    "E8" = 5.718
    "F8" = 7.516
    "G8" = 6.6
    "H8" = 7.631
    "I8" = 6.74
    "J8" = 5.896
    "K8" = 6.246
    "L8" = 6.589
    "M8" = 7.252
    "N8" = 7.882
    "O8" = 8.48
    "B9" = 6.65
    "C9" = 7.15
    "D9" = 7.649
    "E9" = 7.896
    "F9" = 10.38
    "G9" = 9.114000000000001
    "H9" = 10.54
    "I9" = 9.308999999999999
    "J9" = 8.141
    "K9" = 8.625999999999999
    "L9" = 9.1
    "M9" = 10.015
    "N9" = 10.886
    "O9" = 11.714
    "B10" = 4.258
    "C10" = 4.578
    "D10" = 4.897
    "E10" = 5.056
    "F10" = 6.645
    "G10" = 5.835
    "H10" = 6.747
    "I10" = 5.96
    "J10" = 5.213
    "K10" = 5.522
    "L10" = 5.826
    "M10" = 6.411
    "N10" = 6.968
    "O10" = 7.498
    "B11" = 5.656
    "C11" = 6.081
    "D11" = 6.505
    "E11" = 6.715
    "F11" = 8.827
    "G11" = 7.751
    "H11" = 8.962999999999999
    "I11" = 7.917
    "J11" = 6.924
    "K11" = 7.336
    "L11" = 7.739
    "M11" = 8.516999999999999
    "N11" = 9.257
    "O11" = 9.961
    "B12" = 1.238
    "C12" = 1.332
    "D12" = 1.424
    "E12" = 1.47
    "F12" = 1.932
    "G12" = 1.697
    "H12" = 1.962
    "I12" = 1.733
    "J12" = 1.516
    "K12" = 1.606
    "L12" = 1.694
    "M12" = 1.864
    "N12" = 2.026
    "O12" = 2.179
}

foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}
